# Update diagnostic UUID values in the "diagnostics" worksheet (column H)
# to reflect the newly-generated identifiers from the SQLa extensions
# autoloading refactor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$updates = @{
    "H6"  = "1c249b90-02ed-42e3-bc80-8bebe0b3f63b"
    "H7"  = "41d3ee39-c1a6-4232-9e19-3c397bc28e67"
    "H8"  = "65619ffc-d9fe-4b65-83d5-ef481cef90d2"
    "H9"  = "c5290207-8a38-422a-828a-ae6f20805f6c"
    "H10" = "740b7b5d-db7c-461f-8980-59e870a14963"
    "H11" = "130d8033-1979-4035-a8ca-828d9d071b1f"
    "H12" = "42705ffe-03a5-4409-a55a-29eb3faa6365"
    "H13" = "904da14c-df53-4397-82b2-edd4155d5d7a"
    "H14" = "d82c4a03-df84-4d87-958c-9b54baddeb2c"
    "H15" = "5fb4f10b-717a-4aa7-a937-a554d4e2b6f1"
    "H16" = "7bddee20-67e2-4c3d-9d3b-adc300497d57"
    "H17" = "b300f309-cb97-4c1e-9046-42c65feaed8b"

    "H19" = "39ad14b9-3de2-4dea-9075-8182a9f7be7e"
    "H20" = "b4d9c654-79c0-474e-8077-88549890405c"
    "H21" = "4b81795d-4e50-4f8b-b9d4-1d879522eef0"
    "H22" = "17063d26-909d-4cf6-9063-d6a1e0dda382"
    "H23" = "2531d9ad-1bf1-4987-bedd-04baa40f9246"
    "H24" = "51438f74-bab0-40d2-9810-83c888938960"
    "H25" = "c1be0df4-8755-4667-8766-3aa8f61737a5"
    "H26" = "6b54595c-8fc5-493c-bac5-078956ac90e3"
    "H27" = "fa0cc4d3-ea6d-4b2b-b507-48f1c541c42b"
    "H28" = "76434d6d-fb85-401c-9383-9a6c18fa2c48"
    "H29" = "d2702766-09b6-4fed-8a0c-62be3bf50fef"
    "H30" = "0f710509-bb46-4f6b-a48e-4a508a8c3fea"
    "H31" = "439bf9e2-d6c0-41ad-a18b-fcd4590aafa9"
    "H32" = "c9b61112-d008-40bb-8ef0-51a1a60793fc"
    "H33" = "0f830419-702b-46ce-902f-f8ce101713a0"
    "H34" = "cfde60f3-6ef4-4fe2-877d-8b34778687a4"
    "H35" = "07168f2d-4495-48f0-98c1-01c3e1ebaab0"
    "H36" = "70c74ded-5d86-468d-a3db-218adaa673c3"
    "H37" = "98ee1b1f-4b3d-405a-8cd2-db5df12212bf"
    "H38" = "576d7063-64a0-42a4-ac5f-83110f349002"
    "H39" = "1558afa6-3b96-4263-ba04-9585fba3091b"

    "H40" = "621f07f6-3879-4616-bf35-ec0593be5d10"
    "H41" = "bc082cab-267c-4b62-b13b-21446c9cb0c6"
    "H42" = "d77dd9b6-b1a0-4aac-a063-a3da778d82c9"
    "H43" = "c87bd9c4-815e-4878-a0d7-e043e70279ed"
    "H44" = "52ef1e71-a854-41a1-a20f-3ac2796927af"
    "H45" = "04801f92-7dae-47a6-820b-0b1e8524dd26"
    "H46" = "27dfe84b-85cd-4872-b8e6-ef705da3dddc"
    "H47" = "123cbf97-0400-4776-a15d-e181a9830a62"
    "H48" = "e5685493-d5e0-4dec-ae56-721d09fcddaf"
    "H49" = "2620e43a-0b92-46d8-b59e-7c2974a50a17"
    "H50" = "c7f44052-0361-429f-8e45-9023a9278099"
    "H51" = "cc9b2c17-185a-48a7-a7c7-0a432ae8e014"
    "H52" = "ac74606d-0e85-4b55-853e-50bca27c880e"
    "H53" = "83a12958-d579-4d25-8ed2-f8b1fa94b51c"
    "H54" = "d6a90ca6-379d-4844-984a-f84a0a631155"
    "H55" = "97e20099-328f-499a-a988-25d98e06f3f0"
    "H56" = "0cab18ab-0987-456e-a6a8-7346eaa87760"
    "H57" = "727e7674-fc47-44bc-9973-a68bf55d9cc8"
    "H58" = "854d1ffb-566a-4614-b9f8-f29306fce2fc"
    "H59" = "ed5eb2cb-ea5e-48dc-af61-f3baa15fa267"
    "H60" = "85173e0b-5c11-4d69-8349-e876776f4ef7"
    "H61" = "c01e2c32-e028-46d1-a0e7-7f520819fa60"
    "H62" = "6059ded4-437d-4a44-aa93-40d808e2f127"
    "H63" = "d8d44904-32b9-40d5-95d4-b05349628e30"
    "H64" = "558640af-39e9-48e3-87cf-907a383dfa0b"
    "H65" = "8fe53d94-1e80-42ae-8593-28b012c54a3e"
    "H66" = "51c7349b-81cf-4b9a-be85-9f1fbe2203a2"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
